# Update market-price / profit columns (H-N) on several Leve rows across
# the crafting-class sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to
# reflect refreshed Universalis price data pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 326.53845
$ws.Range("J18").Value = 271
$ws.Range("L18").Value = 271
$ws.Range("N18").Value = -839
$ws.Range("H70").Value = 1819.7097
$ws.Range("I70").Value = 1633.3334
$ws.Range("J70").Value = 2077.7693
$ws.Range("K70").Value = 4900.0002
$ws.Range("L70").Value = 6233.3079
$ws.Range("M70").Value = -4630.0002
$ws.Range("N70").Value = -6773.3079
$ws.Range("H73").Value = 1819.7097
$ws.Range("I73").Value = 1633.3334
$ws.Range("J73").Value = 2077.7693
$ws.Range("K73").Value = 4900.0002
$ws.Range("L73").Value = 6233.3079
$ws.Range("M73").Value = -3964.0002
$ws.Range("N73").Value = -8105.3079
$ws.Range("H132").Value = 760.2222
$ws.Range("I132").Value = 708.0606
$ws.Range("K132").Value = 2124.1818
$ws.Range("M132").Value = 405.8181999999997
$ws.Range("H133").Value = 64560.96
$ws.Range("J133").Value = 64560.96
$ws.Range("L133").Value = 64560.96
$ws.Range("N133").Value = -74680.95999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3041.1538
$ws.Range("I2").Value = 1032.2222
$ws.Range("J2").Value = 7561.25
$ws.Range("K2").Value = 1032.2222
$ws.Range("L2").Value = 7561.25
$ws.Range("M2").Value = -919.2221999999999
$ws.Range("N2").Value = -7787.25
$ws.Range("H5").Value = 268.14285
$ws.Range("I5").Value = 300.33334
$ws.Range("J5").Value = 75
$ws.Range("K5").Value = 300.33334
$ws.Range("L5").Value = 75
$ws.Range("M5").Value = -188.33334
$ws.Range("N5").Value = -299
$ws.Range("H116").Value = 3041.1538
$ws.Range("I116").Value = 1032.2222
$ws.Range("J116").Value = 7561.25
$ws.Range("K116").Value = 1032.2222
$ws.Range("L116").Value = 7561.25
$ws.Range("M116").Value = 1261.7778
$ws.Range("N116").Value = -12149.25
$ws.Range("H132").Value = 2190.3333
$ws.Range("I132").Value = 2094.7742
$ws.Range("K132").Value = 6284.3226
$ws.Range("M132").Value = -3754.3226

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3041.1538
$ws.Range("I3").Value = 1032.2222
$ws.Range("J3").Value = 7561.25
$ws.Range("K3").Value = 1032.2222
$ws.Range("L3").Value = 7561.25
$ws.Range("M3").Value = -918.2221999999999
$ws.Range("N3").Value = -7789.25
$ws.Range("H4").Value = 268.14285
$ws.Range("I4").Value = 300.33334
$ws.Range("J4").Value = 75
$ws.Range("K4").Value = 300.33334
$ws.Range("L4").Value = 75
$ws.Range("M4").Value = -185.33334
$ws.Range("N4").Value = -305
$ws.Range("H110").Value = 39702
$ws.Range("J110").Value = 39702
$ws.Range("L110").Value = 39702
$ws.Range("N110").Value = -47882

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("H58").Value = 1822119.5
$ws.Range("I58").Value = 3497586
$ws.Range("J58").Value = 7030.9165
$ws.Range("K58").Value = 3497586
$ws.Range("L58").Value = 7030.9165
$ws.Range("M58").Value = -3497383
$ws.Range("N58").Value = -7436.9165
$ws.Range("H94").Value = 574.63635
$ws.Range("I94").Value = 477.5
$ws.Range("J94").Value = 630.1429
$ws.Range("K94").Value = 477.5
$ws.Range("L94").Value = 630.1429
$ws.Range("M94").Value = -26.5
$ws.Range("N94").Value = -1532.1429
$ws.Range("H99").Value = 3434.8572
$ws.Range("I99").Value = 1903.5
$ws.Range("J99").Value = 5476.6665
$ws.Range("K99").Value = 1903.5
$ws.Range("L99").Value = 5476.6665
$ws.Range("M99").Value = -405.5
$ws.Range("N99").Value = -8472.6665
$ws.Range("H126").Value = 3434.8572
$ws.Range("I126").Value = 1903.5
$ws.Range("J126").Value = 5476.6665
$ws.Range("K126").Value = 5710.5
$ws.Range("L126").Value = 16429.9995
$ws.Range("M126").Value = -3240.5
$ws.Range("N126").Value = -21369.9995
$ws.Range("H132").Value = 3674.08
$ws.Range("I132").Value = 3717.6
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 11152.8
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -8622.8
$ws.Range("N132").Value = -15560
$ws.Range("H136").Value = 1822119.5
$ws.Range("I136").Value = 3497586
$ws.Range("J136").Value = 7030.9165
$ws.Range("K136").Value = 10492758
$ws.Range("L136").Value = 21092.7495
$ws.Range("M136").Value = -10490208
$ws.Range("N136").Value = -26192.7495

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H17").Value = 1243.6364
$ws.Range("I17").Value = 700
$ws.Range("J17").Value = 1896
$ws.Range("K17").Value = 2100
$ws.Range("L17").Value = 5688
$ws.Range("M17").Value = -1931
$ws.Range("N17").Value = -6026
$ws.Range("H129").Value = 1291.0952
$ws.Range("I129").Value = 1306
$ws.Range("J129").Value = 1286.4375
$ws.Range("K129").Value = 3918
$ws.Range("L129").Value = 3859.3125
$ws.Range("M129").Value = 1082
$ws.Range("N129").Value = -13859.3125
$ws.Range("H131").Value = 1105.0741
$ws.Range("I131").Value = 564.44446
$ws.Range("J131").Value = 1213.2
$ws.Range("K131").Value = 1693.33338
$ws.Range("L131").Value = 3639.6
$ws.Range("M131").Value = 3346.66662
$ws.Range("N131").Value = -13719.6
$ws.Range("H132").Value = 2789.1892
$ws.Range("I132").Value = 2500
$ws.Range("J132").Value = 2845.1614
$ws.Range("K132").Value = 22500
$ws.Range("L132").Value = 25606.4526
$ws.Range("M132").Value = -19970
$ws.Range("N132").Value = -30666.4526
$ws.Range("H141").Value = 2849.625
$ws.Range("I141").Value = 2402.4285
$ws.Range("K141").Value = 7207.2855
$ws.Range("M141").Value = -2027.2855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2221.0527
$ws.Range("I113").Value = 2042.8572
$ws.Range("J113").Value = 2720
$ws.Range("K113").Value = 2042.8572
$ws.Range("L113").Value = 2720
$ws.Range("M113").Value = 127.1428000000001
$ws.Range("N113").Value = -7060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 307926.94
$ws.Range("I55").Value = 666736
$ws.Range("J55").Value = 376.2857
$ws.Range("K55").Value = 666736
$ws.Range("L55").Value = 376.2857
$ws.Range("M55").Value = -666563
$ws.Range("N55").Value = -722.2857
$ws.Range("H82").Value = 1600.75
$ws.Range("I82").Value = 700
$ws.Range("J82").Value = 2501.5
$ws.Range("K82").Value = 700
$ws.Range("L82").Value = 2501.5
$ws.Range("M82").Value = -339
$ws.Range("N82").Value = -3223.5
$ws.Range("H85").Value = 1600.75
$ws.Range("I85").Value = 700
$ws.Range("J85").Value = 2501.5
$ws.Range("K85").Value = 700
$ws.Range("L85").Value = 2501.5
$ws.Range("M85").Value = 548
$ws.Range("N85").Value = -4997.5
$ws.Range("H132").Value = 2633.6316
$ws.Range("I132").Value = 1896
$ws.Range("J132").Value = 3647.875
$ws.Range("K132").Value = 5688
$ws.Range("L132").Value = 10943.625
$ws.Range("M132").Value = -3158
$ws.Range("N132").Value = -16003.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 11669150
$ws.Range("I2").Value = 35001450
$ws.Range("J2").Value = 3000
$ws.Range("K2").Value = 35001450
$ws.Range("L2").Value = 3000
$ws.Range("M2").Value = -35001338
$ws.Range("N2").Value = -3224
$ws.Range("H107").Value = 3146.7646
$ws.Range("I107").Value = 1142.8572
$ws.Range("J107").Value = 4549.5
$ws.Range("K107").Value = 3428.5716
$ws.Range("L107").Value = 13648.5
$ws.Range("M107").Value = -1508.5716
$ws.Range("N107").Value = -17488.5
$ws.Range("H132").Value = 2249.75
$ws.Range("I132").Value = 2206.6785
$ws.Range("J132").Value = 2551.25
$ws.Range("K132").Value = 6620.0355
$ws.Range("L132").Value = 7653.75
$ws.Range("M132").Value = -4090.0355
$ws.Range("N132").Value = -12713.75
$ws.Range("H136").Value = 4489.2266
$ws.Range("I136").Value = 1425.68
$ws.Range("J136").Value = 7224.5356
$ws.Range("K136").Value = 4277.04
$ws.Range("L136").Value = 21673.6068
$ws.Range("M136").Value = -1727.04
$ws.Range("N136").Value = -26773.6068
